# Insert a new weekly data row before the existing row 203, shifting
# all subsequent rows (old 203..283) down by one (to 204..284).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new week's record.
$ws.Cells.Item(203, 1).Value2 = 7
$ws.Cells.Item(203, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(203, 3).Value2 = "Ñuble"
$ws.Cells.Item(203, 4).Value2 = 44960
$ws.Cells.Item(203, 5).Value2 = 16
$ws.Cells.Item(203, 6).Value2 = 100112017
$ws.Cells.Item(203, 7).Value2 = "Apio"
$ws.Cells.Item(203, 8).Value2 = "Americana (o)"
$ws.Cells.Item(203, 9).Value2 = "Primera"
$ws.Cells.Item(203, 10).Value2 = 60
$ws.Cells.Item(203, 11).Value2 = 10000
$ws.Cells.Item(203, 12).Value2 = 10000
$ws.Cells.Item(203, 13).Value2 = 10000
$ws.Cells.Item(203, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(203, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(203, 16).Value2 = 1667
$ws.Cells.Item(203, 17).Value2 = 6
$ws.Cells.Item(203, 18).Value2 = "Hortaliza"

# Make sure the new date cell carries the same date number format as
# the rest of column D (style already copied down by the row insert,
# but set it explicitly to be safe).
$ws.Cells.Item(203, 4).NumberFormat = $ws.Cells.Item(204, 4).NumberFormat
